# Corrección a Diebold Mariano y revisión de Cap1
# Updates the dm_completo_General workbook: recomputed DM-test results
# (Matriz_Resultados), p-values (P_valores) and statistics (Estadisticos_DM),
# plus the re-sorted summary (Resumen).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: Matriz_Resultados
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Matriz_Resultados")

$matriz = @{
    "B2"=0;  "C2"=0;  "D2"=0;  "E2"=0;  "F2"=0;  "G2"=0;  "H2"=0;  "I2"=-1; "J2"=0;
    "B3"=0;  "C3"=0;  "D3"=1;  "E3"=1;  "F3"=0;  "G3"=0;  "H3"=0;  "I3"=0;  "J3"=1;
    "B4"=0;  "C4"=-1; "D4"=0;  "E4"=0;  "F4"=0;  "G4"=0;  "H4"=0;  "I4"=0;  "J4"=1;
    "B5"=0;  "C5"=-1; "D5"=0;  "E5"=0;  "F5"=0;  "G5"=0;  "H5"=0;  "I5"=0;  "J5"=0;
    "B6"=0;  "C6"=0;  "D6"=0;  "E6"=0;  "F6"=0;  "G6"=0;  "H6"=0;  "I6"=-1; "J6"=0;
    "B7"=0;  "C7"=0;  "D7"=0;  "E7"=0;  "F7"=0;  "G7"=0;  "H7"=0;  "I7"=0;  "J7"=0;
    "B8"=0;  "C8"=0;  "D8"=0;  "E8"=0;  "F8"=0;  "G8"=0;  "H8"=0;  "I8"=0;  "J8"=0;
    "B9"=1;  "C9"=0;  "D9"=0;  "E9"=0;  "F9"=1;  "G9"=0;  "H9"=0;  "I9"=0;  "J9"=0;
    "B10"=0; "C10"=-1;"D10"=-1;"E10"=0; "F10"=0; "G10"=0; "H10"=0; "I10"=0; "J10"=0;
}
foreach ($addr in $matriz.Keys) {
    $ws1.Range($addr).Value = $matriz[$addr]
}

# ---------------------------------------------------------------------
# Sheet 2: P_valores
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("P_valores")

$pvalores = @{
    "B2"=1;                     "C2"=0.002497270883133229;  "D2"=0.003517808575465464;
    "E2"=0.009285831233703634;  "F2"=0.02069770820294337;   "G2"=0.003563255051977299;
    "H2"=0.003082679151582557;  "I2"=0.00138293586019711;   "J2"=0.02491511647008204;

    "B3"=0.002497270883133229;  "C3"=1;                     "D3"=0.00001922530852516147;
    "E3"=0.000706389124963902;  "F3"=0.002133996810657068;  "G3"=0.0015263875234659;
    "H3"=0.002042881261158591;  "I3"=0.03349818600747323;   "J3"=0.0005746434786069177;

    "B4"=0.003517808575465464;  "C4"=0.00001922530852516147;"D4"=1;
    "E4"=0.001661702955437772;  "F4"=0.003100119261938206;  "G4"=0.007316875113141164;
    "H4"=0.009083781838047811;  "I4"=0.09810747093879968;   "J4"=0.0009178871552482093;

    "B5"=0.009285831233703634;  "C5"=0.000706389124963902;  "D5"=0.001661702955437772;
    "E5"=1;                     "F5"=0.009285844034439217;  "G5"=0.2203322782471007;
    "H5"=0.3456620571406912;    "I5"=0.4936356513560316;    "J5"=0.003547348020978225;

    "B6"=0.02069770820294337;   "C6"=0.002133996810657068;  "D6"=0.003100119261938206;
    "E6"=0.009285844034439217;  "F6"=1;                     "G6"=0.002929803593646474;
    "H6"=0.002468041391740616;  "I6"=0.001044479719741531;  "J6"=0.03116425617557783;

    "B7"=0.003563255051977299;  "C7"=0.0015263875234659;    "D7"=0.007316875113141164;
    "E7"=0.2203322782471007;    "F7"=0.002929803593646474;  "G7"=1;
    "H7"=0.3581911069807318;    "I7"=0.9976026072268833;    "J7"=0.001468590531056524;

    "B8"=0.003082679151582557;  "C8"=0.002042881261158591;  "D8"=0.009083781838047811;
    "E8"=0.3456620571406912;    "F8"=0.002468041391740616;  "G8"=0.3581911069807318;
    "H8"=1;                     "I8"=0.8506648404732422;    "J8"=0.001591974741000701;

    "B9"=0.00138293586019711;   "C9"=0.03349818600747323;   "D9"=0.09810747093879968;
    "E9"=0.4936356513560316;    "F9"=0.001044479719741531;  "G9"=0.9976026072268833;
    "H9"=0.8506648404732422;    "I9"=1;                     "J9"=0.005028328456109454;

    "B10"=0.02491511647008204;  "C10"=0.0005746434786069177;"D10"=0.0009178871552482093;
    "E10"=0.003547348020978225; "F10"=0.03116425617557783;  "G10"=0.001468590531056524;
    "H10"=0.001591974741000701; "I10"=0.005028328456109454; "J10"=1;
}
foreach ($addr in $pvalores.Keys) {
    $ws2.Range($addr).Value = $pvalores[$addr]
}

# ---------------------------------------------------------------------
# Sheet 3: Estadisticos_DM
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Estadisticos_DM")

$estad = @{
    "B2"=0;                    "C2"=3.412238887782399;   "D2"=3.268149501053552;
    "E2"=2.851260737422985;    "F2"=2.492623769941244;   "G2"=3.262725900749976;
    "H2"=3.323827123400642;    "I2"=3.65821535884673;    "J2"=2.407052874157959;

    "B3"=-3.412238887782399;   "C3"=0;                   "D3"=-5.418896875268254;
    "E3"=-3.935136791290341;   "F3"=-3.477949337777977;  "G3"=-3.617322761867777;
    "H3"=-3.496149596096676;   "I3"=-2.267937905842347;  "J3"=-4.019856187300427;

    "B4"=-3.268149501053552;   "C4"=5.418896875268254;   "D4"=0;
    "E4"=-3.582075954620635;   "F4"=-3.321452297101157;  "G4"=-2.955060485683298;
    "H4"=-2.860889571405292;   "I4"=-1.727418275554835;  "J4"=-3.827429307356105;

    "B5"=-2.851260737422985;   "C5"=3.935136791290341;   "D5"=3.582075954620635;
    "E5"=0;                    "F5"=-2.851260133744803;  "G5"=1.261525223544554;
    "H5"=0.9637172936089529;   "I5"=0.6961320670049859;  "J5"=-3.264616567979964;

    "B6"=-2.492623769941244;   "C6"=3.477949337777977;   "D6"=3.321452297101157;
    "E6"=2.851260133744803;    "F6"=0;                   "G6"=3.345222729600425;
    "H6"=3.417168282322994;    "I6"=3.774188311104176;   "J6"=2.30218648941491;

    "B7"=-3.262725900749976;   "C7"=3.617322761867777;   "D7"=2.955060485683298;
    "E7"=-1.261525223544554;   "F7"=-3.345222729600425;  "G7"=0;
    "H7"=-0.9384732931974057;  "I7"=0.003039018182921341;"J7"=-3.633324043832347;

    "B8"=-3.323827123400642;   "C8"=3.496149596096676;   "D8"=2.860889571405292;
    "E8"=-0.9637172936089529;  "F8"=-3.417168282322994;  "G8"=0.9384732931974057;
    "H8"=0;                    "I8"=0.1904992835163427;  "J8"=-3.599871029005107;

    "B9"=-3.65821535884673;    "C9"=2.267937905842347;   "D9"=1.727418275554835;
    "E9"=-0.6961320670049859;  "F9"=-3.774188311104176;  "G9"=-0.003039018182921341;
    "H9"=-0.1904992835163427;  "I9"=0;                   "J9"=-3.11641089207814;

    "B10"=-2.407052874157959;  "C10"=4.019856187300427;  "D10"=3.827429307356105;
    "E10"=3.264616567979964;   "F10"=-2.30218648941491;  "G10"=3.633324043832347;
    "H10"=3.599871029005107;   "I10"=3.11641089207814;   "J10"=0;
}
foreach ($addr in $estad.Keys) {
    $ws3.Range($addr).Value = $estad[$addr]
}

# ---------------------------------------------------------------------
# Sheet 4: Resumen (rows re-sorted by Tasa_Victoria_% after the DM fix)
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("Resumen")

$resumen = @(
    @{A="Sieve Bootstrap";       B=3; C=0; D=5; E=37.5; F=0.5774576320792454},
    @{A="DeepAR";                B=2; C=0; D=6; E=25;   F=1.297661622856752},
    @{A="LSPM";                  B=1; C=1; D=6; E=12.5; F=0.7877186732480409},
    @{A="Block Bootstrapping";   B=0; C=1; D=7; E=0;    F=3.163775868438677},
    @{A="LSPMW";                 B=0; C=1; D=7; E=0;    F=1.489903633896171},
    @{A="AREPD";                 B=0; C=1; D=7; E=0;    F=2.991284383072845},
    @{A="MCPS";                  B=0; C=0; D=8; E=0;    F=1.298246623877678},
    @{A="AV-MCPS";               B=0; C=0; D=8; E=0;    F=1.329647778123854},
    @{A="EnCQR-LSTM";            B=0; C=2; D=6; E=0;    F=2.187906433249951}
)

for ($i = 0; $i -lt $resumen.Count; $i++) {
    $row = $i + 2
    $rec = $resumen[$i]
    $ws4.Range("A$row").Value = $rec.A
    $ws4.Range("B$row").Value = $rec.B
    $ws4.Range("C$row").Value = $rec.C
    $ws4.Range("D$row").Value = $rec.D
    $ws4.Range("E$row").Value = $rec.E
    $ws4.Range("F$row").Value = $rec.F
}

Write-Output "Edit applied"
